$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: plain text id -> new plain text id
$ws.Range("A2").Value = "75c44810a32a3d6447df"

# B2: phone number stored as text (leading "+" must be preserved, so force
# the cell to Text format first or Excel will coerce it to a number and
# drop the "+").
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "+74267426016"

# C2: plain text name -> new plain text name
$ws.Range("C2").Value = "Automation User 10"

# D2: was stored as text "100" in the source file; the edit turns it into a
# genuine numeric value.
$ws.Range("D2").Value = 100

# E2: date stored as literal text "YYYY-MM-DD" (not a real date cell), so
# force Text format first or Excel will turn it into a date serial number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2026-02-16"
